$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, date range) ---
$ws.Range("A8").Characters(21, 2).Text = "17"
$ws.Range("C9").Characters(27, 9).Text = "4/24/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/30/2023"

# L15 simple numeric update (style unchanged)
$ws.Range("L15").Value = -80

# Row 16
$ws.Range("F16").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("L16").Value = 71.428571428571
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -85.882352941176
# Row 17
$ws.Range("C17").Value = 2
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -66.666666666666
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 29
$ws.Range("K17").Value = -13.793103448275
$ws.Range("L17").Value = -13.793103448275
$ws.Range("M17").Value = 8.695652173913
$ws.Range("N17").Value = -62.686567164179
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -42.105263157894
$ws.Range("L18").Value = -15.384615384615
$ws.Range("M18").Value = -31.25
$ws.Range("N18").Value = -93.292682926829
# Row 19
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 15
$ws.Range("H19").Value = 36.363636363636
$ws.Range("I19").Value = 56
$ws.Range("J19").Value = 46
$ws.Range("K19").Value = 21.739130434782
$ws.Range("L19").Value = 64.705882352941
$ws.Range("M19").Value = 93.103448275862
$ws.Range("N19").Value = 33.333333333333
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("F16").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 10
$ws.Range("K20").Value = -28.571428571428
$ws.Range("L20").Value = 66.666666666666
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -91.596638655462
# Row 21
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 23
$ws.Range("G21").Value = 33
$ws.Range("H21").Value = -30.303030303030
$ws.Range("I21").Value = 115
$ws.Range("J21").Value = 122
$ws.Range("K21").Value = -5.737704918032
$ws.Range("L21").Value = 22.340425531914
$ws.Range("M21").Value = 23.655913978494
$ws.Range("N21").Value = -76.288659793814
# Row 23
$ws.Range("F16").Copy($ws.Range("C23"))
$ws.Range("E23").Value = -100
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = -12.5
$ws.Range("L23").Value = -12.5
# Row 24
$ws.Range("C24").Value = 16
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = -42.857142857142
$ws.Range("I24").Value = 161
$ws.Range("J24").Value = 164
$ws.Range("K24").Value = -1.829268292682
$ws.Range("L24").Value = 40
$ws.Range("M24").Value = 103.79746835443
# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 15
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 64
$ws.Range("J25").Value = 51
$ws.Range("K25").Value = 25.490196078431
$ws.Range("L25").Value = 20.754716981132
$ws.Range("M25").Value = -37.254901960784
# Row 26
$ws.Range("L26").Value = -60
# Row 27
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 14.285714285714
$ws.Range("L27").Value = 60
